$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "couleur"
$ws.Range("C6").Value = "RAFALE"
$ws.Range("D6").Value = "Design plat ou les couleurs de l'interface utilisateur plats sont très populaires dans la conception web aujourd'hui où audacieuses, des couleurs vives sont utilisés pour créer des interfaces simples, propres"
$ws.Range("E6").Value = "Fermé"
$ws.Range("F6").Value = "P2"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "couleur"
$ws.Range("C7").Value = "RAFALE"
$ws.Range("D7").Value = "Design plat ou les couleurs de l'interface utilisateur plats sont très populaires dans la conception web aujourd'hui où audacieuses, des couleurs vives sont utilisés pour créer des interfaces simples, propres"
$ws.Range("E7").Value = "Ouvert"
$ws.Range("F7").Value = "P3"
